$d = $word.ActiveDocument

# Trim the trailing clause from the ASICs / Keccak mining sentence:
# "... no ASICs will be created for quite some time, until Smartcash
#  reaches a considerable market cap." ->
# "... no ASICs will be created for quite some time."
$d.Content.Find.Execute(
    "for quite some time, until Smartcash reaches a considerable market cap.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for quite some time.", 2
)
